$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete footer rows 654:659 (sample size / source / elaboration notes)
$ws.Range("A654:D659").EntireRow.Delete() | Out-Null

# Update header row and municipality/state name capitalization, and one numeric precision fix
$ws.Range("A1").Value2 = "mx_state"
$ws.Range("B1").Value2 = "mx_municipality"
$ws.Range("C1").Value2 = "n_matriculas"
$ws.Range("D1").Value2 = "pct_matriculas"
$ws.Range("B7").Value2 = "Pabellón De Arteaga"
$ws.Range("B8").Value2 = "Rincón De Romos"
$ws.Range("B12").Value2 = "Playas De Rosarito"
$ws.Range("B59").Value2 = "Coyame Del Sotol"
$ws.Range("B70").Value2 = "Guadalupe Y Calvo"
$ws.Range("B73").Value2 = "Hidalgo Del Parral"
$ws.Range("B97").Value2 = "San Francisco De Borja"
$ws.Range("B98").Value2 = "San Francisco De Conchos"
$ws.Range("B99").Value2 = "San Francisco Del Oro"
$ws.Range("B107").Value2 = "Valle De Zaragoza"
$ws.Range("B120").Value2 = "San Juan De Sabinas"
$ws.Range("A129").Value2 = "Ciudad De México"
$ws.Range("D143").Value2 = 0.009641255605381168
$ws.Range("B157").Value2 = "Nombre De Dios"
$ws.Range("B160").Value2 = "Pánuco De Coronado"
$ws.Range("B167").Value2 = "San Juan De Guadalupe"
$ws.Range("B168").Value2 = "San Juan Del Río"
$ws.Range("B169").Value2 = "San Luis Del Cordero"
$ws.Range("A178").Value2 = "Estado De México"
$ws.Range("B178").Value2 = "Acambay De Ruíz Castañeda"
$ws.Range("B179").Value2 = "Almoloya De Alquisiras"
$ws.Range("B185").Value2 = "Coacalco De Berriozábal"
$ws.Range("B187").Value2 = "Ecatepec De Morelos"
$ws.Range("B188").Value2 = "Ixtapan De La Sal"
$ws.Range("B189").Value2 = "Naucalpan De Juárez"
$ws.Range("B193").Value2 = "San Felipe Del Progreso"
$ws.Range("B198").Value2 = "Tenango Del Valle"
$ws.Range("B201").Value2 = "Tlalnepantla De Baz"
$ws.Range("B204").Value2 = "Villa De Allende"
$ws.Range("B205").Value2 = "Villa Del Carbón"
$ws.Range("B212").Value2 = "Apaseo El Grande"
$ws.Range("B217").Value2 = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B221").Value2 = "Jaral Del Progreso"
$ws.Range("B231").Value2 = "San Diego De La Unión"
$ws.Range("B233").Value2 = "San Francisco Del Rincón"
$ws.Range("B234").Value2 = "San Luis De La Paz"
$ws.Range("B236").Value2 = "Santa Cruz De Juventino Rosas"
$ws.Range("B237").Value2 = "Silao De La Victoria"
$ws.Range("B241").Value2 = "Valle De Santiago"
$ws.Range("B244").Value2 = "Acapulco De Juárez"
$ws.Range("B245").Value2 = "Ajuchitlán Del Progreso"
$ws.Range("B247").Value2 = "Atoyac De Álvarez"
$ws.Range("B248").Value2 = "Ayutla De Los Libres"
$ws.Range("B249").Value2 = "Chilapa De Álvarez"
$ws.Range("B250").Value2 = "Chilpancingo De Los Bravo"
$ws.Range("B253").Value2 = "Coyuca De Catalán"
$ws.Range("B254").Value2 = "Cutzamala De Pinzón"
$ws.Range("B257").Value2 = "Iguala De La Independencia"
$ws.Range("B262").Value2 = "Taxco De Alarcón"
$ws.Range("B270").Value2 = "Cuautepec De Hinojosa"
$ws.Range("B277").Value2 = "Mineral Del Monte"
$ws.Range("B278").Value2 = "Mixquiahuala De Juárez"
$ws.Range("B279").Value2 = "Pachuca De Soto"
$ws.Range("B283").Value2 = "Tulancingo De Bravo"
$ws.Range("B285").Value2 = "Zacualtipán De Ángeles"
$ws.Range("B289").Value2 = "Atotonilco El Alto"
$ws.Range("B290").Value2 = "Autlán De Navarro"
$ws.Range("B298").Value2 = "Encarnación De Díaz"
$ws.Range("B301").Value2 = "Huejuquilla El Alto"
$ws.Range("B302").Value2 = "Ixtlahuacán De Los Membrillos"
$ws.Range("B303").Value2 = "Ixtlahuacán Del Río"
$ws.Range("B309").Value2 = "Lagos De Moreno"
$ws.Range("B316").Value2 = "San Miguel El Alto"
$ws.Range("B317").Value2 = "Tamazula De Gordiano"
$ws.Range("B319").Value2 = "Teocuitatlán De Corona"
$ws.Range("B320").Value2 = "Tepatitlán De Morelos"
$ws.Range("B323").Value2 = "Tizapán El Alto"
$ws.Range("B324").Value2 = "Tlajomulco De Zúñiga"
$ws.Range("B328").Value2 = "Unión De San Antonio"
$ws.Range("B329").Value2 = "Unión De Tula"
$ws.Range("B334").Value2 = "Yahualica De González Gallo"
$ws.Range("B336").Value2 = "Zapotlán El Grande"
$ws.Range("B345").Value2 = "Coalcomán De Vázquez Pallares"
$ws.Range("B346").Value2 = "Cojumatlán De Régules"
$ws.Range("B402").Value2 = "San Nicolás De Los Garza"
$ws.Range("B406").Value2 = "El Barrio De La Soledad"
$ws.Range("B408").Value2 = "Guevea De Humboldt"
$ws.Range("B409").Value2 = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B410").Value2 = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B412").Value2 = "Mariscala De Juárez"
$ws.Range("B413").Value2 = "Miahuatlán De Porfirio Díaz"
$ws.Range("B414").Value2 = "Oaxaca De Juárez"
$ws.Range("B415").Value2 = "Ocotlán De Morelos"
$ws.Range("B442").Value2 = "Tlacolula De Matamoros"
$ws.Range("B444").Value2 = "Villa De Etla"
$ws.Range("B445").Value2 = "Villa Tejúpam De La Unión"
$ws.Range("B453").Value2 = "Chila De La Sal"
$ws.Range("B457").Value2 = "Huehuetlán El Chico"
$ws.Range("B458").Value2 = "Huehuetlán El Grande"
$ws.Range("B461").Value2 = "Los Reyes De Juárez"
$ws.Range("B472").Value2 = "Tepexi De Rodríguez"
$ws.Range("B473").Value2 = "Tetela De Ocampo"
$ws.Range("B485").Value2 = "Amealco De Bonfil"
$ws.Range("B487").Value2 = "Jalpan De Serra"
$ws.Range("B488").Value2 = "Landa De Matamoros"
$ws.Range("B492").Value2 = "San Juan Del Río"
$ws.Range("B499").Value2 = "Ciudad Del Maíz"
$ws.Range("B505").Value2 = "Mexquitic De Carmona"
$ws.Range("B511").Value2 = "San Ciro De Acosta"
$ws.Range("B515").Value2 = "Villa De Guadalupe"
$ws.Range("B516").Value2 = "Villa De Ramos"
$ws.Range("B547").Value2 = "Nacozari De García"
$ws.Range("B553").Value2 = "San Pedro De La Cueva"
$ws.Range("B575").Value2 = "Tepetitla De Lardizábal"
$ws.Range("B582").Value2 = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B585").Value2 = "Amatlán De Los Reyes"
$ws.Range("B589").Value2 = "Cosamaloapan De Carpio"
$ws.Range("B596").Value2 = "Hueyapan De Ocampo"
$ws.Range("B599").Value2 = "Lerdo De Tejada"
$ws.Range("B600").Value2 = "Martínez De La Torre"
$ws.Range("B609").Value2 = "Sayula De Alemán"
$ws.Range("B621").Value2 = "Cañitas De Felipe Pescador"
$ws.Range("B624").Value2 = "El Plateado De Joaquín Amaro"
$ws.Range("B636").Value2 = "Moyahua De Estrada"
$ws.Range("B637").Value2 = "Nochistlán De Mejía"
$ws.Range("B638").Value2 = "Noria De Ángeles"
$ws.Range("B648").Value2 = "Villa De Cos"
